$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose IDs (RM 232, SC 92) were dropped from the source data
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-apply the full target grid of values (including the new missing-value mask)
$ws.Range("B2").Value = -19.7
$ws.Range("C2").Value = 14.9
$ws.Range("D2").Value = -13.5
$ws.Range("E2").Value = -7.2
$ws.Range("F2").Value = 18.03
$ws.Range("B3").Value = -19.7
$ws.Range("C3").Value = 11.2
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 17.64
$ws.Range("B4").Value = -18.7
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = -15.4
$ws.Range("E4").Value = -6.4
$ws.Range("F4").Value = 17.97
$ws.Range("B5").Value = -19.5
$ws.Range("C5").Value = 12.3
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = -5
$ws.Range("F5").Value = 17.66
$ws.Range("B6").Value = -19.8
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("B7").Value = -19.5
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = -13.8
$ws.Range("E7").Value = -7.1
$ws.Range("F7").Value = 17.24
$ws.Range("B8").Value = -19.9
$ws.Range("C8").Value = 15.5
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").Value = 17.05
$ws.Range("B9").Value = -20.6
$ws.Range("C9").Value = 10.5
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("B10").Value = -19.8
$ws.Range("C10").Value = 11.5
$ws.Range("D10").Value = -14.7
$ws.Range("E10").Value = -6.1
$ws.Range("F10").Value = 16.43
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 11.4
$ws.Range("D11").Value = -15.5
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = 17.65
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("D12").Value = -14.1
$ws.Range("E12").Value = -5.3
$ws.Range("F12").Value = 17.45
$ws.Range("B13").Value = -19.9
$ws.Range("C13").Value = 12.5
$ws.Range("D13").Value = -13.9
$ws.Range("E13").Value = -5.3
$ws.Range("F13").Value = 17.1
$ws.Range("B14").Value = -19.6
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = -13.1
$ws.Range("E14").Value = -5.4
$ws.Range("F14").Value = 17.76
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 12.5
$ws.Range("D15").Value = -15.2
$ws.Range("E15").Value = -8.4
$ws.Range("F15").Value = 16.2
$ws.Range("B16").Value = -19.5
$ws.Range("C16").Value = 13.5
$ws.Range("D16").Value = -15.3
$ws.Range("E16").Value = -6.9
$ws.Range("F16").Value = 17.34
$ws.Range("B17").Value = -19.4
$ws.Range("C17").Value = 11.2
$ws.Range("D17").Value = -14.7
$ws.Range("E17").Value = -7.3
$ws.Range("F17").Value = 17.78
$ws.Range("B18").Value = -19.6
$ws.Range("C18").Value = 11.5
$ws.Range("D18").Value = -15.2
$ws.Range("E18").Value = -8.5
$ws.Range("F18").Value = 18.35
$ws.Range("B19").Value = -20.6
$ws.Range("C19").Value = 13.2
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").Value = 17.81
$ws.Range("B20").Value = -19
$ws.Range("C20").Value = 12.5
$ws.Range("D20").Value = -14
$ws.Range("E20").Value = -7.2
$ws.Range("F20").Value = 17.73
$ws.Range("B21").Value = -18.9
$ws.Range("C21").Value = 12.7
$ws.Range("D21").Value = -14.3
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F21").Value = 16.58
$ws.Range("B22").Value = -19.3
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("E22").Value = -6.1
$ws.Range("F22").Value = 16.81
$ws.Range("B23").Value = -19.5
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48
$ws.Range("B24").Value = -17.7
$ws.Range("C24").ClearContents()
$ws.Range("D24").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("F24").Value = 16.78
$ws.Range("B25").Value = -19.5
$ws.Range("C25").Value = 10.7
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38
$ws.Range("B27").ClearContents()
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("D29").ClearContents()
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").ClearContents()
$ws.Range("F31").Value = 17.18
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
